$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 1.02, 1.061690715476676, 1.0391534462867, 1.066598522113135, 1.07572469191559, 1.042476319329654, 1.066665092628052, 1.04193944350433, 1.06930863974435, 1.078410505242127, 1.06817987940778),
    @(3, 1.02, 1.063258837415286, 1.039771940143554, 1.068037772727517, 1.077346774571287, 1.04280809889427, 1.067884407854698, 1.042368603338219, 1.07056195920964, 1.079847987716948, 1.069400926201906),
    @(4, 1.02, 1.064271812929896, 1.040171699114588, 1.068967753785769, 1.078395295376777, 1.043020993392603, 1.068671287369704, 1.042645159589666, 1.071371109082895, 1.080776571786889, 1.070188923176079),
    @(5, 1.02, 1.06469726884863, 1.040339651364668, 1.069358411602882, 1.0788358446818, 1.043110067721289, 1.069001595889526, 1.042761152411787, 1.071710844088027, 1.081166582090045, 1.070519700771886),
    @(6, 1.02, 1.064768681633466, 1.040367845069014, 1.069423987028833, 1.078909800503891, 1.043124998733428, 1.069057027278321, 1.042780612235565, 1.071767861980464, 1.081232045273516, 1.070575210879611),
    @(7, 1.02, 1.064277499449327, 1.040173943717975, 1.068972974968646, 1.078401182985367, 1.043022185281872, 1.068675702907327, 1.042646710556469, 1.071375650326059, 1.080781784553884, 1.070193344984271),
    @(8, 1.02, 1.062221027940126, 1.039362562418401, 1.067085199422392, 1.076273110118431, 1.042588817305223, 1.067077604670951, 1.042084716626524, 1.069732588338552, 1.078896637555024, 1.068592977265096),
    @(9, 1.02, 1.058583804920159, 1.03792933684284, 1.063748325902971, 1.072514555214098, 1.041811384711141, 1.064245166719951, 1.04108563369406, 1.066822964971808, 1.07556244471911, 1.06575651692732),
    @(10, 1.02, 1.056149373014146, 1.036971465667563, 1.061516297615146, 1.070002517615583, 1.041283715136071, 1.062345424159928, 1.040413596896595, 1.064873139964258, 1.073330871711025, 1.063854076514958),
    @(11, 1.02, 1.055092829077323, 1.03655611723499, 1.060547932814812, 1.068913154448613, 1.041052977388749, 1.061520001373543, 1.04012116039477, 1.064026356639538, 1.072362383733097, 1.06302748153354),
    @(12, 1.02, 1.054700007179917, 1.036401749363512, 1.060187947089768, 1.068508260450292, 1.040966930281426, 1.061212970190951, 1.040012318431953, 1.063711440082188, 1.072002303749232, 1.062720014331452),
    @(13, 1.02, 1.05478428593251, 1.036434865842375, 1.060265178608405, 1.068595123341723, 1.040985403138621, 1.06127884912609, 1.040035675295489, 1.063779008245384, 1.072079557687811, 1.062785986822234),
    @(14, 1.02, 1.055060365999275, 1.036543358960863, 1.060518182291094, 1.068879691067941, 1.041045871674297, 1.061494630950494, 1.040112167941521, 1.06400033341193, 1.07233262640962, 1.06300207508158),
    @(15, 1.02, 1.055230418092143, 1.036610193324666, 1.060674027251157, 1.069054988396667, 1.041083083101901, 1.061627523701156, 1.040159268644723, 1.064136648075119, 1.072488505003582, 1.063135156555191),
    @(16, 1.02, 1.05621944045362, 1.036999018567395, 1.061580524366472, 1.070074779735894, 1.041298980769768, 1.062400144579566, 1.040432974460881, 1.0649292848684, 1.073395099856049, 1.063908874643869),
    @(17, 1.02, 1.056839172695744, 1.037242761173346, 1.062148635770201, 1.070714022979498, 1.041433802776534, 1.062884027828075, 1.040604276157921, 1.065425810699845, 1.073963186541415, 1.064393445062093),
    @(18, 1.02, 1.057200419393115, 1.037384875861012, 1.062479824597643, 1.071086725773836, 1.041512224914127, 1.063165997061816, 1.040704054701628, 1.065715185544985, 1.074294330070956, 1.064675814724488),
    @(19, 1.02, 1.057323556054763, 1.03743332379283, 1.062592721014978, 1.071213781485029, 1.041538928054676, 1.063262095457506, 1.040738053106881, 1.065813814390547, 1.074407205816588, 1.064772049590922),
    @(20, 1.02, 1.05677270543357, 1.037216615725836, 1.062087701575498, 1.070645454541671, 1.041419360139801, 1.062832139897284, 1.040585911479742, 1.06537256309117, 1.073902258195585, 1.064341483444488),
    @(21, 1.02, 1.054979077726057, 1.036511412907759, 1.060443687145233, 1.068795900117061, 1.041028074624928, 1.061431100586433, 1.040089648795822, 1.063935169322608, 1.072258113453412, 1.062938454497115),
    @(22, 1.02, 1.05384918089186, 1.03606750851334, 1.059408336439001, 1.067631525368774, 1.040780084774927, 1.060547705534841, 1.039776366481814, 1.063029199576351, 1.071222399988458, 1.062053804923223),
    @(23, 1.02, 1.054448370159551, 1.03630287990859, 1.059957359061084, 1.068248926874207, 1.040911736687478, 1.061016250357042, 1.039942563573167, 1.063509684872553, 1.071771641662761, 1.062523015132814),
    @(24, 1.02, 1.056802739847257, 1.037228429902239, 1.062115235674487, 1.070676438152453, 1.041425886817072, 1.062855586659334, 1.040594210116483, 1.065396624136091, 1.073929789747119, 1.06436496350363),
    @(25, 1.02, 1.05952576729556, 1.03830027663266, 1.064612265220047, 1.073487308402219, 1.042014014711963, 1.064979406002591, 1.041344968933802, 1.067576916167084, 1.076425924289583, 1.06649179891396)
)

$cols = @(2,3,4,5,6,9,10,11,12,13,14)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $c = $cols[$i]
        $v = $entry[$i + 1]
        $ws.Cells.Item($r, $c).Value = $v
    }
}
